$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 611 (rows 611-710 shift down to 614-713).
$ws.Rows.Item(611).Resize(3).Insert()

# Common (constant across this data block) column values A-J, K.
$mercadoId = 8
$mercado = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$codreg = 4
$tipo = "Fruta"
$productoId = 100106
$producto = "Oleaginosos"
$categoriaId = 100106002
$categoria = "Palta"
$variedad = "Hass"

# New weekly data rows inserted at 611-613.
$newRows = @(
  @{ Row = 611; Fecha = 44637; Calidad = "Especial"; Volumen = 240; PMin = 3300; PMax = 3400; PProm = 3350; Unidad = "`$/kilo (en caja de 17 kilos)"; Origen = "Provincia de Limarí"; PKg = 3350; KgUnidad = 1 },
  @{ Row = 612; Fecha = 44637; Calidad = "Primera";  Volumen = 300; PMin = 3100; PMax = 3200; PProm = 3150; Unidad = "`$/kilo (en caja de 17 kilos)"; Origen = "Provincia de Limarí"; PKg = 3150; KgUnidad = 1 },
  @{ Row = 613; Fecha = 44637; Calidad = "Segunda";  Volumen = 240; PMin = 2900; PMax = 3000; PProm = 2950; Unidad = "`$/kilo (en caja de 17 kilos)"; Origen = "Provincia de Limarí"; PKg = 2950; KgUnidad = 1 }
)

foreach ($r in $newRows) {
  $row = $r.Row
  $ws.Cells.Item($row, 1).Value = $mercadoId
  $ws.Cells.Item($row, 2).Value = $mercado
  $ws.Cells.Item($row, 3).Value = $region
  $ws.Cells.Item($row, 4).Value = $r.Fecha
  $ws.Cells.Item($row, 5).Value = $codreg
  $ws.Cells.Item($row, 6).Value = $tipo
  $ws.Cells.Item($row, 7).Value = $productoId
  $ws.Cells.Item($row, 8).Value = $producto
  $ws.Cells.Item($row, 9).Value = $categoriaId
  $ws.Cells.Item($row, 10).Value = $categoria
  $ws.Cells.Item($row, 11).Value = $variedad
  $ws.Cells.Item($row, 12).Value = $r.Calidad
  $ws.Cells.Item($row, 13).Value = $r.Volumen
  $ws.Cells.Item($row, 14).Value = $r.PMin
  $ws.Cells.Item($row, 15).Value = $r.PMax
  $ws.Cells.Item($row, 16).Value = $r.PProm
  $ws.Cells.Item($row, 17).Value = $r.Unidad
  $ws.Cells.Item($row, 18).Value = $r.Origen
  $ws.Cells.Item($row, 19).Value = $r.PKg
  $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
